$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with the new values
# A2 must remain stored as text (it was an inline string before), so force
# the cell's number format to Text before assigning the value.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "12"
$ws.Range("B2").Value = "Thiago"
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = 7
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = 8
$ws.Range("G2").Value = 7.5

# Remove row 3 entirely (shifts rows up, deleting the row)
$ws.Rows("3:3").Delete()
